$wb = $excel.ActiveWorkbook

# Rename the two sheets
$wsMonth = $wb.Worksheets.Item("Month_Wise_Concentration")
$wsMonth.Name = "Sales_Month_Wise_Concentration"

$wsPlant = $wb.Worksheets.Item("Plant_Wise_Concentration")
$wsPlant.Name = "Sales_Plant_Wise_Concentration"

# Update the config sheet (Main) cells that mirror the sheet names
$wsMain = $wb.Worksheets.Item("Main")
$wsMain.Range("B133").Value = "Sales_Month_Wise_Concentration"
$wsMain.Range("B134").Value = "Sales_Plant_Wise_Concentration"
